$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '29.124.11'
$ws.Range("E2").Value2 = '  +0.03%  '
$ws.Range("D3").Value2 = '1.831.10'
$ws.Range("E3").Value2 = '  -0.45%  '
$ws.Range("D4").Value2 = '''1.001'
$ws.Range("E4").Value2 = '  +0.01%  '
$ws.Range("D5").Value2 = '''243.31'
$ws.Range("E5").Value2 = '  -0.07%  '
$ws.Range("D6").Value2 = '''0.6257'
$ws.Range("E6").Value2 = '  -0.54%  '
$ws.Range("D7").Value2 = '''1.004'
$ws.Range("E7").Value2 = '  +0.11%  '
$ws.Range("D8").Value2 = '''0.07454'
$ws.Range("E8").Value2 = '  -1.75%  '
$ws.Range("E9").Value2 = '  -0.35%  '
$ws.Range("D10").Value2 = '''23.28'
$ws.Range("E10").Value2 = '  +2.93%  '
$ws.Range("D11").Value2 = '''0.07706'
$ws.Range("E11").Value2 = '  -0.72%  '
$ws.Range("D12").Value2 = '1.838.80'
$ws.Range("E12").Value2 = '  -0.18%  '
$ws.Range("D13").Value2 = '''5.007'
$ws.Range("E13").Value2 = '  +0.64%  '
$ws.Range("D14").Value2 = '''0.6673'
$ws.Range("D15").Value2 = '''82.55'
$ws.Range("E15").Value2 = '  -0.91%  '
$ws.Range("D16").Value2 = '''0.000009351'
$ws.Range("E16").Value2 = '  -6.54%  '
$ws.Range("D17").Value2 = '''5.957'
$ws.Range("E17").Value2 = '  -1.91%  '
$ws.Range("D18").Value2 = '29.126.45'
$ws.Range("E18").Value2 = '  +0.01%  '
$ws.Range("D19").Value2 = '2.072.71'
$ws.Range("E19").Value2 = '  -0.63%  '
$ws.Range("D20").Value2 = '''12.58'
$ws.Range("E20").Value2 = '  +1.16%  '
$ws.Range("D21").Value2 = '''222.91'
$ws.Range("E21").Value2 = '  -1.86%  '
$ws.Range("E22").Value2 = '  +0.07%  '
$ws.Range("D23").Value2 = '''7.125'
$ws.Range("E23").Value2 = '  -1.39%  '
$ws.Range("E24").Value2 = '  +0.02%  '
$ws.Range("D25").Value2 = '''160.16'
$ws.Range("E25").Value2 = '  +0.09%  '
$ws.Range("D26").Value2 = '''0.1391'
$ws.Range("E26").Value2 = '  +0.28%  '
$ws.Range("D27").Value2 = '''8.492'
$ws.Range("E27").Value2 = '  -0.30%  '
$ws.Range("E28").Value2 = '  -0.48%  '
$ws.Range("D29").Value2 = '''1.493'
$ws.Range("E29").Value2 = '  -0.17%  '
$ws.Range("D30").Value2 = '''0.05815'
$ws.Range("E30").Value2 = '  +10.40%  '
$ws.Range("D31").Value2 = '''4.156'
$ws.Range("E31").Value2 = '  +1.18%  '
$ws.Range("D32").Value2 = '''4.125'
$ws.Range("E32").Value2 = '  +2.60%  '
$ws.Range("E33").Value2 = '  +1.31%  '
$ws.Range("D34").Value2 = '''0.7387'
$ws.Range("E34").Value2 = '  +0.16%  '
$ws.Range("D35").Value2 = '''1.826'
$ws.Range("E35").Value2 = '  -1.16%  '
$ws.Range("E36").Value2 = '  -0.12%  '
$ws.Range("D37").Value2 = '''2.676'
$ws.Range("E37").Value2 = '  -0.28%  '
$ws.Range("B38").Value2 = 'MXToken'
$ws.Range("C38").Value2 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value2 = '''2.765'
$ws.Range("E38").Value2 = '  -0.07%  '
$ws.Range("B39").Value2 = 'Maker'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value2 = '1.225.59'
$ws.Range("E39").Value2 = '  -1.61%  '
$ws.Range("D40").Value2 = '''0.01771'
$ws.Range("E40").Value2 = '  -0.97%  '
$ws.Range("D41").Value2 = '''6.491'
$ws.Range("E41").Value2 = '  +2.03%  '
$ws.Range("D42").Value2 = '''0.8932'
$ws.Range("E42").Value2 = '  -0.93%  '
$ws.Range("D43").Value2 = '''1.003'
$ws.Range("E43").Value2 = '  +0.11%  '
$ws.Range("B44").Value2 = 'BabyDogeCoin'
$ws.Range("C44").Value2 = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").Value2 = '''0.00000000130'
$ws.Range("E44").Value2 = '  +3.36%  '
$ws.Range("B45").Value2 = 'Quant'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value2 = '''102.13'
$ws.Range("E45").Value2 = '  +0.02%  '
$ws.Range("D46").Value2 = '1.978.32'
$ws.Range("E46").Value2 = '  -0.40%  '
$ws.Range("B47").Value2 = 'Aave'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value2 = '''65.89'
$ws.Range("E47").Value2 = '  +2.43%  '
$ws.Range("D48").Value2 = '''0.5091'
$ws.Range("E48").Value2 = '  -0.76%  '
$ws.Range("D49").Value2 = '''0.07591'
$ws.Range("E49").Value2 = '  +13.76%  '
$ws.Range("D50").Value2 = '''0.4060'
$ws.Range("E50").Value2 = '  +0.28%  '
$ws.Range("D51").Value2 = '''8.988'
$ws.Range("E51").Value2 = '  +0.68%  '
